# Append a freshly-scraped Lancers listing at the top of the job list
# (row 8, just under the unchanged rows 2-7) and stamp every row with the
# new scrape timestamp, matching the "Append: 2025-11-19 01:19 JST" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-11-19 01:19:49"

# 1) Make room for the new listing: push the old rows 8-14 down to 9-15.
$ws.Rows.Item(8).Insert()

# 2) Fill in the brand-new listing in row 8.
$ws.Range("A8").Value = $newTimestamp
$ws.Range("B8").Value = "【急募】Android用のライブ壁紙アプリ開発エンジニアを探しています!"
$ws.Range("C8").Value = "システム開発"
$ws.Range("D8").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E8").Value = "期限情報なし"
$ws.Range("F8").Value = "https://www.lancers.jp/work/detail/5436829"
$ws.Range("G8").Value = 93
$ws.Range("H8").Value = "◆開発 ◇アプリ"

# 3) Refresh the scrape timestamp on every other data row (2-7, now 9-15).
$ws.Range("A2").Value = $newTimestamp
$ws.Range("A3").Value = $newTimestamp
$ws.Range("A4").Value = $newTimestamp
$ws.Range("A5").Value = $newTimestamp
$ws.Range("A6").Value = $newTimestamp
$ws.Range("A7").Value = $newTimestamp
$ws.Range("A9").Value = $newTimestamp
$ws.Range("A10").Value = $newTimestamp
$ws.Range("A11").Value = $newTimestamp
$ws.Range("A12").Value = $newTimestamp
$ws.Range("A13").Value = $newTimestamp
$ws.Range("A14").Value = $newTimestamp
$ws.Range("A15").Value = $newTimestamp

# 4) Rebuild the hyperlinks on column F (the row-insert above does not
#    re-anchor the existing hyperlink objects to their shifted rows, so
#    drop the stale collection and re-create one link per data row).
$ws.Range("F2").Hyperlinks.Delete()

$links = @(
    @{ Row = 2;  Url = "https://www.lancers.jp/work/detail/5436501" },
    @{ Row = 3;  Url = "https://www.lancers.jp/work/detail/5436391" },
    @{ Row = 4;  Url = "https://www.lancers.jp/work/detail/5436668" },
    @{ Row = 5;  Url = "https://www.lancers.jp/work/detail/5429809" },
    @{ Row = 6;  Url = "https://www.lancers.jp/work/detail/5436149" },
    @{ Row = 7;  Url = "https://www.lancers.jp/work/detail/5436594" },
    @{ Row = 8;  Url = "https://www.lancers.jp/work/detail/5436829" },
    @{ Row = 9;  Url = "https://www.lancers.jp/work/detail/5425629" },
    @{ Row = 10; Url = "https://www.lancers.jp/work/detail/5436021" },
    @{ Row = 11; Url = "https://www.lancers.jp/work/detail/5341051" },
    @{ Row = 12; Url = "https://www.lancers.jp/work/detail/5436476" },
    @{ Row = 13; Url = "https://www.lancers.jp/work/detail/5436366" },
    @{ Row = 14; Url = "https://www.lancers.jp/work/detail/5436426" },
    @{ Row = 15; Url = "https://www.lancers.jp/work/detail/5436248" }
)

foreach ($link in $links) {
    $cell = $ws.Cells.Item($link.Row, 6)
    $ws.Hyperlinks.Add($cell, $link.Url) | Out-Null
    $cell.Style = "Hyperlink"
}
